$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text, even if it looks like a number,
# by forcing a Text number format for the duration of the write and then
# restoring the cell to the default (Normal) style so formatting is unaffected.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$ws.Cells.Item(2, 4).Value = "29.222.04"
$ws.Cells.Item(2, 5).Value = "  -0.89%  "
$ws.Cells.Item(3, 4).Value = "1.867.15"
$ws.Cells.Item(3, 5).Value = "  -0.44%  "
Set-TextValue $ws.Cells.Item(4, 4) "0.9999"
$ws.Cells.Item(4, 5).Value = "  +0.15%  "
Set-TextValue $ws.Cells.Item(5, 4) "0.7099"
$ws.Cells.Item(5, 5).Value = "  -0.79%  "
Set-TextValue $ws.Cells.Item(6, 4) "241.47"
$ws.Cells.Item(6, 5).Value = "  -0.10%  "
$ws.Cells.Item(7, 5).Value = "  +0.11%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3115"
$ws.Cells.Item(8, 5).Value = "  +0.24%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.07665"
$ws.Cells.Item(9, 5).Value = "  -3.76%  "
Set-TextValue $ws.Cells.Item(10, 4) "24.68"
$ws.Cells.Item(10, 5).Value = "  -2.65%  "
Set-TextValue $ws.Cells.Item(11, 4) "0.08368"
$ws.Cells.Item(11, 5).Value = "  +1.23%  "
$ws.Cells.Item(12, 4).Value = "1.860.04"
$ws.Cells.Item(12, 5).Value = "  -0.62%  "
Set-TextValue $ws.Cells.Item(13, 4) "5.231"
$ws.Cells.Item(14, 5).Value = "  -2.69%  "
Set-TextValue $ws.Cells.Item(15, 4) "91.18"
$ws.Cells.Item(15, 5).Value = "  +0.03%  "
$ws.Cells.Item(16, 4).Value = "29.230.71"
$ws.Cells.Item(16, 5).Value = "  -0.78%  "
Set-TextValue $ws.Cells.Item(17, 4) "5.943"
$ws.Cells.Item(17, 5).Value = "  +0.21%  "
Set-TextValue $ws.Cells.Item(18, 4) "243.62"
$ws.Cells.Item(18, 5).Value = "  -0.92%  "
Set-TextValue $ws.Cells.Item(19, 4) "0.000007833"
$ws.Cells.Item(20, 4).Value = "2.115.19"
$ws.Cells.Item(20, 5).Value = "  +0.52%  "
Set-TextValue $ws.Cells.Item(21, 4) "13.11"
$ws.Cells.Item(21, 5).Value = "  -1.85%  "
Set-TextValue $ws.Cells.Item(22, 4) "0.9994"
$ws.Cells.Item(22, 5).Value = "  +0.09%  "
Set-TextValue $ws.Cells.Item(23, 4) "7.853"
$ws.Cells.Item(23, 5).Value = "  -1.91%  "
$ws.Cells.Item(24, 5).Value = "  +0.22%  "
Set-TextValue $ws.Cells.Item(25, 4) "0.1637"
$ws.Cells.Item(25, 5).Value = "  +1.22%  "
Set-TextValue $ws.Cells.Item(26, 4) "163.30"
$ws.Cells.Item(26, 5).Value = "  -0.09%  "
Set-TextValue $ws.Cells.Item(27, 4) "8.961"
Set-TextValue $ws.Cells.Item(28, 4) "18.52"
$ws.Cells.Item(28, 5).Value = "  +1.11%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.508"
$ws.Cells.Item(29, 5).Value = "  +1.06%  "
Set-TextValue $ws.Cells.Item(30, 4) "1.310"
$ws.Cells.Item(30, 5).Value = "  -3.52%  "
Set-TextValue $ws.Cells.Item(31, 4) "4.394"
$ws.Cells.Item(31, 5).Value = "  +0.20%  "
$ws.Cells.Item(32, 5).Value = "  +3.05%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.05133"
$ws.Cells.Item(33, 5).Value = "  -2.46%  "
Set-TextValue $ws.Cells.Item(34, 4) "0.7937"
$ws.Cells.Item(34, 5).Value = "  +9.10%  "
Set-TextValue $ws.Cells.Item(35, 4) "1.912"
$ws.Cells.Item(35, 5).Value = "  -2.10%  "
Set-TextValue $ws.Cells.Item(36, 4) "1.166"
Set-TextValue $ws.Cells.Item(37, 4) "2.686"
$ws.Cells.Item(37, 5).Value = "  +0.43%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.01855"
$ws.Cells.Item(38, 5).Value = "  -0.86%  "
Set-TextValue $ws.Cells.Item(39, 4) "2.709"
$ws.Cells.Item(39, 5).Value = "  +0.04%  "
$ws.Cells.Item(40, 4).Value = "1.153.21"
$ws.Cells.Item(40, 5).Value = "  -5.45%  "
Set-TextValue $ws.Cells.Item(41, 4) "6.379"
$ws.Cells.Item(41, 5).Value = "  +4.24%  "
Set-TextValue $ws.Cells.Item(42, 4) "0.8955"
$ws.Cells.Item(42, 5).Value = "  -1.72%  "
Set-TextValue $ws.Cells.Item(43, 4) "73.16"
$ws.Cells.Item(43, 5).Value = "  -0.76%  "
Set-TextValue $ws.Cells.Item(44, 4) "0.9994"
$ws.Cells.Item(44, 5).Value = "  +0.04%  "
Set-TextValue $ws.Cells.Item(45, 4) "103.24"
$ws.Cells.Item(45, 5).Value = "  +0.98%  "
$ws.Cells.Item(46, 4).Value = "2.012.38"
$ws.Cells.Item(46, 5).Value = "  +0.10%  "
Set-TextValue $ws.Cells.Item(47, 4) "0.5166"
$ws.Cells.Item(48, 5).Value = "  -1.03%  "
Set-TextValue $ws.Cells.Item(49, 4) "9.342"
$ws.Cells.Item(49, 5).Value = "  +0.03%  "
$ws.Cells.Item(50, 5).Value = "  +2.32%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.4294"
